$wb = $excel.ActiveWorkbook

$wsBoM = $wb.Worksheets.Item("BoM")
$wsDNF = $wb.Worksheets.Item("DNF")

# BoM sheet - "Net Name" column (X) updates
$wsBoM.Range("X15").Value = "/RESET2,/MISO2,+5V,/MOSI2,GND,/SCK2"
$wsBoM.Range("X16").Value = "Net-(J3-Pin_5),Net-(J3-Pin_1),Net-(J3-Pin_3),Net-(J3-Pin_2),Net-(J3-Pin_4)"
$wsBoM.Range("X17").Value = "Net-(J6-Pin_3),Net-(J6-Pin_2),Net-(J6-Pin_4),Net-(J6-Pin_6),Net-(J6-Pin_1),Net-(J6-Pin_5)"
$wsBoM.Range("X21").Value = "Net-(U1-UCAP),Net-(J3-Pin_5),/RESET2,VBUS,Net-(J3-Pin_2),Net-(J4-Pin_1),Net-(J3-Pin_4),Net-(J4-Pin_2),Net-(U1-D+),/SCK2,/MISO2,/MOSI2,/RXLED,/TXLED,Net-(J6-Pin_3),unconnected-(U1-PB0-Pad14),Net-(U1-XTAL1),Net-(J3-Pin_3),Net-(J6-Pin_2),Net-(J4-Pin_3),Net-(J6-Pin_4),Net-(U1-D-),Net-(J6-Pin_6),GND,Net-(J6-Pin_5),Earth,/DTR,Net-(J3-Pin_1),+5V,Net-(J4-Pin_4),Net-(U1-PC0{slash}XTAL2)"

# DNF sheet - "Net Name" column (X) updates
$wsDNF.Range("X10").Value = "VBUS,Net-(J2-VBUS)"
$wsDNF.Range("X12").Value = "Net-(J4-Pin_2),Net-(J4-Pin_1),Net-(J4-Pin_3),Net-(J4-Pin_4)"
$wsDNF.Range("X13").Value = "Net-(J2-D-),Net-(J2-Shield),Net-(J2-D+),Net-(J2-VBUS),Earth"
$wsDNF.Range("X14").Value = "Net-(U1-D+),Net-(J2-D+)"
$wsDNF.Range("X15").Value = "Net-(U1-PC0{slash}XTAL2),Net-(U1-XTAL1)"
$wsDNF.Range("X16").Value = "Net-(J2-Shield),Net-(J2-D+)"
$wsDNF.Range("X17").Value = "Net-(U1-PC0{slash}XTAL2),Net-(U1-XTAL1)"
